# ----------------------------------------------------------------------------
# Adds 22 new "opgave" rows (130-151) covering several new exam topics, plus
# the corresponding new shared-string entries, mirroring the existing table
# layout/styling exactly. Also refreshes the active sheet view/selection.
# ----------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template cells: reuse the existing (already-present) cell styles so no new
# style entries are minted -- PasteSpecial(formats) copies the exact style index.
$styleSrc3 = $ws.Cells.Item(2, 4)
$styleSrc4 = $ws.Cells.Item(2, 6)
$styleSrc5 = $ws.Cells.Item(2, 9)
$styleSrc6 = $ws.Cells.Item(2, 11)
$styleSrc7 = $ws.Cells.Item(22, 11)

# Row 130
$ws.Cells.Item(130, 1).Value = 'Akkerranden'
$ws.Cells.Item(130, 2).Value = 1
$ws.Cells.Item(130, 3).Value = 'Toon aan'
$ws.Cells.Item(130, 4).Value = 'B1: Rekenen
C1: Tabellen'
$styleSrc3.Copy()
$ws.Cells.Item(130, 4).PasteSpecial(-4122)
$ws.Cells.Item(130, 5).Value = 'Tabel
Eenheden
Verhouding'
$styleSrc3.Copy()
$ws.Cells.Item(130, 5).PasteSpecial(-4122)
$ws.Cells.Item(130, 6).Value = 'Oppervlakte berekenen, gegevens uit tekst en tabel combineren en rekenen, uitkomst vergelijken met waarde'
$styleSrc4.Copy()
$ws.Cells.Item(130, 6).PasteSpecial(-4122)
$ws.Cells.Item(130, 7).Value = 'Deel 1 H1'
$ws.Cells.Item(130, 8).Value = '12.1'
$ws.Cells.Item(130, 9).Value = 'Rekenregels en verhoudingen'
$styleSrc5.Copy()
$ws.Cells.Item(130, 9).PasteSpecial(-4122)
$ws.Cells.Item(130, 10).Value = 'Allerlei formules'
$styleSrc5.Copy()
$ws.Cells.Item(130, 10).PasteSpecial(-4122)
$ws.Cells.Item(130, 11).Value = '5 havo 3'
$styleSrc6.Copy()
$ws.Cells.Item(130, 11).PasteSpecial(-4122)
$ws.Cells.Item(130, 12).Value = 'H2'
$styleSrc6.Copy()
$ws.Cells.Item(130, 12).PasteSpecial(-4122)
$styleSrc5.Copy()
$ws.Cells.Item(130, 14).PasteSpecial(-4122)
$ws.Rows.Item(130).RowHeight = 80

# Row 131
$ws.Cells.Item(131, 2).Value = 2
$ws.Cells.Item(131, 3).Value = 'Bereken'
$ws.Cells.Item(131, 4).Value = 'C3: Formules en vergelijkingen
C4: Lineaire verbanden'
$styleSrc3.Copy()
$ws.Cells.Item(131, 4).PasteSpecial(-4122)
$ws.Cells.Item(131, 5).Value = 'Lineair verband
Meer dan twee variabelen'
$styleSrc3.Copy()
$ws.Cells.Item(131, 5).PasteSpecial(-4122)
$ws.Cells.Item(131, 6).Value = 'Waarden invullen in formule'
$styleSrc4.Copy()
$ws.Cells.Item(131, 6).PasteSpecial(-4122)
$ws.Cells.Item(131, 7).Value = 'Deel 3 H11'
$ws.Cells.Item(131, 8).Value = '12.3'
$ws.Cells.Item(131, 9).Value = 'Formules en variabelen'
$styleSrc5.Copy()
$ws.Cells.Item(131, 9).PasteSpecial(-4122)
$ws.Cells.Item(131, 10).Value = 'Allerlei formules'
$styleSrc5.Copy()
$ws.Cells.Item(131, 10).PasteSpecial(-4122)
$ws.Cells.Item(131, 11).Value = '5 havo 3'
$styleSrc6.Copy()
$ws.Cells.Item(131, 11).PasteSpecial(-4122)
$ws.Cells.Item(131, 12).Value = 'H2'
$styleSrc6.Copy()
$ws.Cells.Item(131, 12).PasteSpecial(-4122)
$styleSrc5.Copy()
$ws.Cells.Item(131, 14).PasteSpecial(-4122)
$ws.Rows.Item(131).RowHeight = 34

# Row 132
$ws.Cells.Item(132, 2).Value = 3
$ws.Cells.Item(132, 3).Value = 'Bereken'
$ws.Cells.Item(132, 4).Value = 'C3: Formules en vergelijkingen
C4: Lineaire verbanden'
$styleSrc3.Copy()
$ws.Cells.Item(132, 4).PasteSpecial(-4122)
$ws.Cells.Item(132, 5).Value = 'Lineair verband
Meer dan twee variabelen
Grafische rekenmachine
Vergelijking
Afronden'
$styleSrc3.Copy()
$ws.Cells.Item(132, 5).PasteSpecial(-4122)
$ws.Cells.Item(132, 6).Value = 'Waarde invullen in formule, vergelijking oplossen (GR)'
$styleSrc4.Copy()
$ws.Cells.Item(132, 6).PasteSpecial(-4122)
$ws.Cells.Item(132, 7).Value = 'Deel 1 H3'
$ws.Cells.Item(132, 8).Value = '12.3'
$ws.Cells.Item(132, 9).Value = 'Tabellen en grafieken'
$styleSrc5.Copy()
$ws.Cells.Item(132, 9).PasteSpecial(-4122)
$ws.Cells.Item(132, 10).Value = 'Allerlei formules'
$styleSrc5.Copy()
$ws.Cells.Item(132, 10).PasteSpecial(-4122)
$ws.Cells.Item(132, 11).Value = '5 havo 3'
$styleSrc6.Copy()
$ws.Cells.Item(132, 11).PasteSpecial(-4122)
$ws.Cells.Item(132, 12).Value = 'H2'
$styleSrc6.Copy()
$ws.Cells.Item(132, 12).PasteSpecial(-4122)
$styleSrc5.Copy()
$ws.Cells.Item(132, 14).PasteSpecial(-4122)
$ws.Rows.Item(132).RowHeight = 85

# Row 133
$ws.Cells.Item(133, 2).Value = 4
$ws.Cells.Item(133, 3).Value = 'Herleid'
$ws.Cells.Item(133, 4).Value = 'B2: Algebra
C3: Formules en vergelijkingen
C4: Lineaire verbanden'
$styleSrc3.Copy()
$ws.Cells.Item(133, 4).PasteSpecial(-4122)
$ws.Cells.Item(133, 5).Value = 'Lineair verband
Meer dan twee variabelen
Herleiden'
$styleSrc3.Copy()
$ws.Cells.Item(133, 5).PasteSpecial(-4122)
$ws.Cells.Item(133, 6).Value = 'Waarde invullen in formule, variabele vrij maken in lineaire vergelijking'
$styleSrc4.Copy()
$ws.Cells.Item(133, 6).PasteSpecial(-4122)
$ws.Cells.Item(133, 7).Value = 'Deel 2 H5'
$ws.Cells.Item(133, 8).Value = '12.2'
$ws.Cells.Item(133, 9).Value = 'Lineaire verbanden'
$styleSrc5.Copy()
$ws.Cells.Item(133, 9).PasteSpecial(-4122)
$ws.Cells.Item(133, 10).Value = 'Allerlei formules'
$styleSrc5.Copy()
$ws.Cells.Item(133, 10).PasteSpecial(-4122)
$ws.Cells.Item(133, 11).Value = '5 havo 3'
$styleSrc6.Copy()
$ws.Cells.Item(133, 11).PasteSpecial(-4122)
$ws.Cells.Item(133, 12).Value = 'H2'
$styleSrc6.Copy()
$ws.Cells.Item(133, 12).PasteSpecial(-4122)
$styleSrc5.Copy()
$ws.Cells.Item(133, 14).PasteSpecial(-4122)
$ws.Rows.Item(133).RowHeight = 51

# Row 134
$ws.Cells.Item(134, 1).Value = 'Onderzoek naar rekenvaardigheid'
$ws.Cells.Item(134, 2).Value = 5
$ws.Cells.Item(134, 3).Value = 'Beredeneer'
$ws.Cells.Item(134, 4).Value = 'E1: Datapresentaties interpreteren'
$styleSrc3.Copy()
$ws.Cells.Item(134, 4).PasteSpecial(-4122)
$ws.Cells.Item(134, 5).Value = 'Staafdiagram/histogram
Statistische visualisatie
Conclusies
Redeneren
Gemiddelde'
$styleSrc3.Copy()
$ws.Cells.Item(134, 5).PasteSpecial(-4122)
$ws.Cells.Item(134, 6).Value = 'Redeneren over gemiddelde aan de hand van staafdiagram'
$styleSrc4.Copy()
$ws.Cells.Item(134, 6).PasteSpecial(-4122)
$ws.Cells.Item(134, 7).Value = 'Deel 1 H2'
$ws.Cells.Item(134, 8).Value = '12.5'
$ws.Cells.Item(134, 9).Value = 'Verwerken van data'
$styleSrc5.Copy()
$ws.Cells.Item(134, 9).PasteSpecial(-4122)
$ws.Cells.Item(134, 10).Value = 'Verdelingen'
$styleSrc5.Copy()
$ws.Cells.Item(134, 10).PasteSpecial(-4122)
$ws.Cells.Item(134, 11).Value = '5 havo 2'
$styleSrc6.Copy()
$ws.Cells.Item(134, 11).PasteSpecial(-4122)
$ws.Cells.Item(134, 12).Value = 'H7'
$styleSrc6.Copy()
$ws.Cells.Item(134, 12).PasteSpecial(-4122)
$styleSrc5.Copy()
$ws.Cells.Item(134, 14).PasteSpecial(-4122)
$ws.Rows.Item(134).RowHeight = 85

# Row 135
$ws.Cells.Item(135, 2).Value = 6
$ws.Cells.Item(135, 3).Value = 'Beredeneer'
$ws.Cells.Item(135, 4).Value = 'C1: Tabellen
E3: Data en verdelingen'
$styleSrc3.Copy()
$ws.Cells.Item(135, 4).PasteSpecial(-4122)
$ws.Cells.Item(135, 5).Value = 'Tabel
Normale verdeling'
$styleSrc3.Copy()
$ws.Cells.Item(135, 5).PasteSpecial(-4122)
$ws.Cells.Item(135, 6).Value = 'Beredeneren of een variabele normaal verdeeld is aan de hand van tabel'
$styleSrc4.Copy()
$ws.Cells.Item(135, 6).PasteSpecial(-4122)
$ws.Cells.Item(135, 7).Value = 'Deel 2 H6'
$ws.Cells.Item(135, 8).Value = '12.5'
$ws.Cells.Item(135, 9).Value = 'Statistiek en beslissingen'
$styleSrc5.Copy()
$ws.Cells.Item(135, 9).PasteSpecial(-4122)
$ws.Cells.Item(135, 10).Value = 'Verdelingen'
$styleSrc5.Copy()
$ws.Cells.Item(135, 10).PasteSpecial(-4122)
$ws.Cells.Item(135, 11).Value = '5 havo 2'
$styleSrc6.Copy()
$ws.Cells.Item(135, 11).PasteSpecial(-4122)
$ws.Cells.Item(135, 12).Value = 'H7'
$styleSrc6.Copy()
$ws.Cells.Item(135, 12).PasteSpecial(-4122)
$styleSrc5.Copy()
$ws.Cells.Item(135, 14).PasteSpecial(-4122)
$ws.Rows.Item(135).RowHeight = 48

# Row 136
$ws.Cells.Item(136, 2).Value = 7
$ws.Cells.Item(136, 3).Value = 'Bepaal'
$ws.Cells.Item(136, 4).Value = 'E4: Statistische uitspraken'
$styleSrc3.Copy()
$ws.Cells.Item(136, 4).PasteSpecial(-4122)
$ws.Cells.Item(136, 5).Value = 'Tabel
Groepen vergelijken
Boxplots vergelijken
Boxplot
Effectgrootte'
$styleSrc3.Copy()
$ws.Cells.Item(136, 5).PasteSpecial(-4122)
$ws.Cells.Item(136, 6).Value = 'Methodes voor groepen vergelijken kiezen, effectgrootte berekenen, boxplots maken en vergelijken'
$styleSrc4.Copy()
$ws.Cells.Item(136, 6).PasteSpecial(-4122)
$ws.Cells.Item(136, 7).Value = 'Deel 3 H10'
$ws.Cells.Item(136, 8).Value = '12.5'
$ws.Cells.Item(136, 9).Value = 'Statistische variabelen'
$styleSrc5.Copy()
$ws.Cells.Item(136, 9).PasteSpecial(-4122)
$ws.Cells.Item(136, 10).Value = 'Conclusies uit data'
$styleSrc5.Copy()
$ws.Cells.Item(136, 10).PasteSpecial(-4122)
$ws.Cells.Item(136, 11).Value = '5 havo 5'
$styleSrc6.Copy()
$ws.Cells.Item(136, 11).PasteSpecial(-4122)
$ws.Cells.Item(136, 12).Value = 'H7'
$styleSrc6.Copy()
$ws.Cells.Item(136, 12).PasteSpecial(-4122)
$styleSrc5.Copy()
$ws.Cells.Item(136, 14).PasteSpecial(-4122)
$ws.Rows.Item(136).RowHeight = 85

# Row 137
$ws.Cells.Item(137, 2).Value = 8
$ws.Cells.Item(137, 3).Value = 'Beredeneer'
$ws.Cells.Item(137, 4).Value = 'E2: Data verwerken'
$styleSrc3.Copy()
$ws.Cells.Item(137, 4).PasteSpecial(-4122)
$ws.Cells.Item(137, 5).Value = 'Tabel
Spreiding
Standaardafwijking
Kwartielafspraak'
$styleSrc3.Copy()
$ws.Cells.Item(137, 5).PasteSpecial(-4122)
$ws.Cells.Item(137, 6).Value = 'Twee spreidingsmaten kiezen, spreiding vergelijken'
$styleSrc4.Copy()
$ws.Cells.Item(137, 6).PasteSpecial(-4122)
$ws.Cells.Item(137, 7).Value = 'Deel 1 H2'
$ws.Cells.Item(137, 8).Value = '12.5'
$ws.Cells.Item(137, 9).Value = 'Verwerken van data'
$styleSrc5.Copy()
$ws.Cells.Item(137, 9).PasteSpecial(-4122)
$ws.Cells.Item(137, 10).Value = 'Verdelingen'
$styleSrc5.Copy()
$ws.Cells.Item(137, 10).PasteSpecial(-4122)
$ws.Cells.Item(137, 11).Value = '5 havo 2'
$styleSrc6.Copy()
$ws.Cells.Item(137, 11).PasteSpecial(-4122)
$ws.Cells.Item(137, 12).Value = 'H7'
$styleSrc6.Copy()
$ws.Cells.Item(137, 12).PasteSpecial(-4122)
$styleSrc5.Copy()
$ws.Cells.Item(137, 14).PasteSpecial(-4122)
$ws.Rows.Item(137).RowHeight = 68

# Row 138
$ws.Cells.Item(138, 2).Value = 9
$ws.Cells.Item(138, 3).Value = 'Beredeneer'
$ws.Cells.Item(138, 4).Value = 'E2: Data verwerken'
$styleSrc3.Copy()
$ws.Cells.Item(138, 4).PasteSpecial(-4122)
$ws.Cells.Item(138, 5).Value = 'Grafiek
Spreiding'
$styleSrc3.Copy()
$ws.Cells.Item(138, 5).PasteSpecial(-4122)
$ws.Cells.Item(138, 6).Value = 'Spreidingsmaat kiezen, spreiding van twee scores vergelijken'
$styleSrc4.Copy()
$ws.Cells.Item(138, 6).PasteSpecial(-4122)
$ws.Cells.Item(138, 7).Value = 'Deel 1 H2'
$ws.Cells.Item(138, 8).Value = '12.5'
$ws.Cells.Item(138, 9).Value = 'Verwerken van data'
$styleSrc5.Copy()
$ws.Cells.Item(138, 9).PasteSpecial(-4122)
$ws.Cells.Item(138, 10).Value = 'Verdelingen'
$styleSrc5.Copy()
$ws.Cells.Item(138, 10).PasteSpecial(-4122)
$ws.Cells.Item(138, 11).Value = '5 havo 2'
$styleSrc6.Copy()
$ws.Cells.Item(138, 11).PasteSpecial(-4122)
$ws.Cells.Item(138, 12).Value = 'H7'
$styleSrc6.Copy()
$ws.Cells.Item(138, 12).PasteSpecial(-4122)
$styleSrc5.Copy()
$ws.Cells.Item(138, 14).PasteSpecial(-4122)
$ws.Rows.Item(138).RowHeight = 48

# Row 139
$ws.Cells.Item(139, 1).Value = 'Great Barrier Reef'
$ws.Cells.Item(139, 2).Value = 10
$ws.Cells.Item(139, 3).Value = 'Bereken'
$ws.Cells.Item(139, 4).Value = 'B1: Rekenen'
$styleSrc3.Copy()
$ws.Cells.Item(139, 4).PasteSpecial(-4122)
$ws.Cells.Item(139, 5).Value = 'Procentuele verandering
Afronden'
$styleSrc3.Copy()
$ws.Cells.Item(139, 5).PasteSpecial(-4122)
$ws.Cells.Item(139, 6).Value = 'Gegevens uit tekst verwerken, procentuele verandering berekenen'
$styleSrc4.Copy()
$ws.Cells.Item(139, 6).PasteSpecial(-4122)
$ws.Cells.Item(139, 7).Value = 'Deel 1 H3'
$ws.Cells.Item(139, 8).Value = '12.1'
$ws.Cells.Item(139, 9).Value = 'Tabellen en grafieken'
$styleSrc5.Copy()
$ws.Cells.Item(139, 9).PasteSpecial(-4122)
$ws.Cells.Item(139, 10).Value = 'Lineaire en exponentiele formules'
$styleSrc5.Copy()
$ws.Cells.Item(139, 10).PasteSpecial(-4122)
$ws.Cells.Item(139, 11).Value = '5 havo 1'
$styleSrc6.Copy()
$ws.Cells.Item(139, 11).PasteSpecial(-4122)
$ws.Cells.Item(139, 12).Value = 'H5'
$styleSrc6.Copy()
$ws.Cells.Item(139, 12).PasteSpecial(-4122)
$styleSrc5.Copy()
$ws.Cells.Item(139, 14).PasteSpecial(-4122)
$ws.Rows.Item(139).RowHeight = 48

# Row 140
$ws.Cells.Item(140, 2).Value = 11
$ws.Cells.Item(140, 3).Value = 'Bereken'
$styleSrc3.Copy()
$ws.Cells.Item(140, 3).PasteSpecial(-4122)
$ws.Cells.Item(140, 4).Value = 'C5: Exponentiële verbanden'
$styleSrc3.Copy()
$ws.Cells.Item(140, 4).PasteSpecial(-4122)
$ws.Cells.Item(140, 5).Value = 'Exponentieel verband
Groeifactor
Afronden'
$styleSrc3.Copy()
$ws.Cells.Item(140, 5).PasteSpecial(-4122)
$ws.Cells.Item(140, 6).Value = 'Groeifactor omrekenen naar andere tijdseenheid, omrekenen naar groeipercentage'
$styleSrc4.Copy()
$ws.Cells.Item(140, 6).PasteSpecial(-4122)
$ws.Cells.Item(140, 7).Value = 'Deel 3 H9'
$styleSrc3.Copy()
$ws.Cells.Item(140, 7).PasteSpecial(-4122)
$ws.Cells.Item(140, 8).Value = '12.4'
$styleSrc3.Copy()
$ws.Cells.Item(140, 8).PasteSpecial(-4122)
$ws.Cells.Item(140, 9).Value = 'Formules en variabelen'
$styleSrc5.Copy()
$ws.Cells.Item(140, 9).PasteSpecial(-4122)
$ws.Cells.Item(140, 10).Value = 'Lineaire en exponentiele formules'
$styleSrc5.Copy()
$ws.Cells.Item(140, 10).PasteSpecial(-4122)
$ws.Cells.Item(140, 11).Value = '5 havo 1'
$styleSrc6.Copy()
$ws.Cells.Item(140, 11).PasteSpecial(-4122)
$ws.Cells.Item(140, 12).Value = 'H5'
$styleSrc6.Copy()
$ws.Cells.Item(140, 12).PasteSpecial(-4122)
$styleSrc5.Copy()
$ws.Cells.Item(140, 14).PasteSpecial(-4122)
$ws.Rows.Item(140).RowHeight = 51

# Row 141
$ws.Cells.Item(141, 2).Value = 12
$ws.Cells.Item(141, 3).Value = 'Bereken'
$ws.Cells.Item(141, 4).Value = 'C3: Formules en vergelijkingen
C5: Exponentiële verbanden'
$styleSrc3.Copy()
$ws.Cells.Item(141, 4).PasteSpecial(-4122)
$ws.Cells.Item(141, 5).Value = 'Groeifactor
Exponentieel verband
Vergelijking
Grafische rekenmachine
Verdubbelingstijd'
$styleSrc3.Copy()
$ws.Cells.Item(141, 5).PasteSpecial(-4122)
$ws.Cells.Item(141, 6).Value = 'Groeifactor bepalen, vergelijking opstellen, oplossen (GR)'
$styleSrc4.Copy()
$ws.Cells.Item(141, 6).PasteSpecial(-4122)
$ws.Cells.Item(141, 7).Value = 'Deel 3 H9'
$ws.Cells.Item(141, 8).Value = '12.4'
$ws.Cells.Item(141, 9).Value = 'Formules en variabelen'
$styleSrc5.Copy()
$ws.Cells.Item(141, 9).PasteSpecial(-4122)
$ws.Cells.Item(141, 10).Value = 'Lineaire en exponentiele formules'
$styleSrc5.Copy()
$ws.Cells.Item(141, 10).PasteSpecial(-4122)
$ws.Cells.Item(141, 11).Value = '5 havo 1'
$styleSrc6.Copy()
$ws.Cells.Item(141, 11).PasteSpecial(-4122)
$ws.Cells.Item(141, 12).Value = 'H5'
$styleSrc6.Copy()
$ws.Cells.Item(141, 12).PasteSpecial(-4122)
$styleSrc5.Copy()
$ws.Cells.Item(141, 14).PasteSpecial(-4122)
$ws.Rows.Item(141).RowHeight = 85

# Row 142
$ws.Cells.Item(142, 1).Value = 'Studieschuld'
$ws.Cells.Item(142, 2).Value = 13
$ws.Cells.Item(142, 3).Value = 'Bereken'
$styleSrc3.Copy()
$ws.Cells.Item(142, 3).PasteSpecial(-4122)
$ws.Cells.Item(142, 4).Value = 'C5: Exponentiële verbanden'
$styleSrc3.Copy()
$ws.Cells.Item(142, 4).PasteSpecial(-4122)
$ws.Cells.Item(142, 5).Value = 'Groeifactor
Afronden'
$styleSrc3.Copy()
$ws.Cells.Item(142, 5).PasteSpecial(-4122)
$ws.Cells.Item(142, 6).Value = 'Groeifactor omrekenen naar andere tijdseenheid'
$styleSrc4.Copy()
$ws.Cells.Item(142, 6).PasteSpecial(-4122)
$ws.Cells.Item(142, 7).Value = 'Deel 3 H9'
$styleSrc3.Copy()
$ws.Cells.Item(142, 7).PasteSpecial(-4122)
$ws.Cells.Item(142, 8).Value = '12.4'
$styleSrc3.Copy()
$ws.Cells.Item(142, 8).PasteSpecial(-4122)
$ws.Cells.Item(142, 9).Value = 'Formules en variabelen'
$styleSrc5.Copy()
$ws.Cells.Item(142, 9).PasteSpecial(-4122)
$ws.Cells.Item(142, 10).Value = 'Lineaire en exponentiele formules'
$styleSrc5.Copy()
$ws.Cells.Item(142, 10).PasteSpecial(-4122)
$ws.Cells.Item(142, 11).Value = '5 havo 1'
$styleSrc6.Copy()
$ws.Cells.Item(142, 11).PasteSpecial(-4122)
$ws.Cells.Item(142, 12).Value = 'H5'
$styleSrc6.Copy()
$ws.Cells.Item(142, 12).PasteSpecial(-4122)
$styleSrc5.Copy()
$ws.Cells.Item(142, 14).PasteSpecial(-4122)
$ws.Rows.Item(142).RowHeight = 34

# Row 143
$ws.Cells.Item(143, 2).Value = 14
$ws.Cells.Item(143, 3).Value = 'Bereken'
$ws.Cells.Item(143, 4).Value = 'C3: Formules en vergelijkingen
C5: Exponentiële verbanden'
$styleSrc3.Copy()
$ws.Cells.Item(143, 4).PasteSpecial(-4122)
$ws.Cells.Item(143, 5).Value = 'Exponentieel verband
Vergelijking
Grafische rekenmachine'
$styleSrc3.Copy()
$ws.Cells.Item(143, 5).PasteSpecial(-4122)
$ws.Cells.Item(143, 6).Value = 'Waarde invullen, vergelijking oplossen (GR)'
$styleSrc4.Copy()
$ws.Cells.Item(143, 6).PasteSpecial(-4122)
$ws.Cells.Item(143, 7).Value = 'Deel 3 H9'
$ws.Cells.Item(143, 8).Value = '12.4'
$ws.Cells.Item(143, 9).Value = 'Formules en variabelen'
$styleSrc5.Copy()
$ws.Cells.Item(143, 9).PasteSpecial(-4122)
$ws.Cells.Item(143, 10).Value = 'Lineaire en exponentiele formules'
$styleSrc5.Copy()
$ws.Cells.Item(143, 10).PasteSpecial(-4122)
$ws.Cells.Item(143, 11).Value = '5 havo 1'
$styleSrc6.Copy()
$ws.Cells.Item(143, 11).PasteSpecial(-4122)
$ws.Cells.Item(143, 12).Value = 'H5'
$styleSrc6.Copy()
$ws.Cells.Item(143, 12).PasteSpecial(-4122)
$styleSrc5.Copy()
$ws.Cells.Item(143, 14).PasteSpecial(-4122)
$ws.Rows.Item(143).RowHeight = 51

# Row 144
$ws.Cells.Item(144, 2).Value = 15
$ws.Cells.Item(144, 3).Value = 'Lees af'
$styleSrc3.Copy()
$ws.Cells.Item(144, 3).PasteSpecial(-4122)
$ws.Cells.Item(144, 4).Value = 'C1: Tabellen'
$styleSrc3.Copy()
$ws.Cells.Item(144, 4).PasteSpecial(-4122)
$ws.Cells.Item(144, 5).Value = 'Tabel'
$styleSrc3.Copy()
$ws.Cells.Item(144, 5).PasteSpecial(-4122)
$ws.Cells.Item(144, 6).Value = 'Waarde in tabel aflezen'
$styleSrc4.Copy()
$ws.Cells.Item(144, 6).PasteSpecial(-4122)
$ws.Cells.Item(144, 7).Value = 'Deel 1 H3'
$styleSrc3.Copy()
$ws.Cells.Item(144, 7).PasteSpecial(-4122)
$ws.Cells.Item(144, 8).Value = '12.1'
$styleSrc3.Copy()
$ws.Cells.Item(144, 8).PasteSpecial(-4122)
$ws.Cells.Item(144, 9).Value = 'Tabellen en grafieken'
$styleSrc5.Copy()
$ws.Cells.Item(144, 9).PasteSpecial(-4122)
$ws.Cells.Item(144, 10).Value = 'Lineaire en exponentiele formules'
$styleSrc5.Copy()
$ws.Cells.Item(144, 10).PasteSpecial(-4122)
$ws.Cells.Item(144, 11).Value = '5 havo 1'
$styleSrc6.Copy()
$ws.Cells.Item(144, 11).PasteSpecial(-4122)
$ws.Cells.Item(144, 12).Value = 'H2'
$styleSrc6.Copy()
$ws.Cells.Item(144, 12).PasteSpecial(-4122)
$styleSrc5.Copy()
$ws.Cells.Item(144, 14).PasteSpecial(-4122)
$ws.Rows.Item(144).RowHeight = 17

# Row 145
$ws.Cells.Item(145, 2).Value = 16
$ws.Cells.Item(145, 3).Value = 'Bereken'
$ws.Cells.Item(145, 4).Value = 'C4: Lineaire verbanden'
$styleSrc3.Copy()
$ws.Cells.Item(145, 4).PasteSpecial(-4122)
$ws.Cells.Item(145, 5).Value = 'Tabel
Lineair verband
Lineaire inter-/extrapoleren
Richtingscoëfficiënt'
$styleSrc3.Copy()
$ws.Cells.Item(145, 5).PasteSpecial(-4122)
$ws.Cells.Item(145, 6).Value = 'Twee waarden uit tabel aflezen, lineair interpoleren'
$styleSrc4.Copy()
$ws.Cells.Item(145, 6).PasteSpecial(-4122)
$ws.Cells.Item(145, 7).Value = 'Deel 2 H5'
$ws.Cells.Item(145, 8).Value = '12.2'
$ws.Cells.Item(145, 9).Value = 'Lineaire verbanden'
$styleSrc5.Copy()
$ws.Cells.Item(145, 9).PasteSpecial(-4122)
$ws.Cells.Item(145, 10).Value = 'Lineaire en exponentiele formules'
$styleSrc5.Copy()
$ws.Cells.Item(145, 10).PasteSpecial(-4122)
$ws.Cells.Item(145, 11).Value = '5 havo 1'
$styleSrc6.Copy()
$ws.Cells.Item(145, 11).PasteSpecial(-4122)
$ws.Cells.Item(145, 12).Value = 'H2'
$styleSrc6.Copy()
$ws.Cells.Item(145, 12).PasteSpecial(-4122)
$styleSrc5.Copy()
$ws.Cells.Item(145, 14).PasteSpecial(-4122)
$ws.Rows.Item(145).RowHeight = 68

# Row 146
$ws.Cells.Item(146, 1).Value = 'Papierformaten'
$ws.Cells.Item(146, 2).Value = 17
$ws.Cells.Item(146, 3).Value = 'Bereken'
$styleSrc3.Copy()
$ws.Cells.Item(146, 3).PasteSpecial(-4122)
$ws.Cells.Item(146, 4).Value = 'C5: Exponentiële verbanden'
$styleSrc3.Copy()
$ws.Cells.Item(146, 4).PasteSpecial(-4122)
$ws.Cells.Item(146, 5).Value = 'Tabel
Exponentieel verband
Groeifactor'
$styleSrc3.Copy()
$ws.Cells.Item(146, 5).PasteSpecial(-4122)
$ws.Cells.Item(146, 6).Value = 'Waarde uit tabel halen, hoeveelheid 11 keer halveren'
$styleSrc4.Copy()
$ws.Cells.Item(146, 6).PasteSpecial(-4122)
$ws.Cells.Item(146, 7).Value = 'Deel 3 H9'
$styleSrc3.Copy()
$ws.Cells.Item(146, 7).PasteSpecial(-4122)
$ws.Cells.Item(146, 8).Value = '12.4'
$styleSrc3.Copy()
$ws.Cells.Item(146, 8).PasteSpecial(-4122)
$ws.Cells.Item(146, 9).Value = 'Formules en variabelen'
$styleSrc5.Copy()
$ws.Cells.Item(146, 9).PasteSpecial(-4122)
$ws.Cells.Item(146, 10).Value = 'Lineaire en exponentiele formules'
$styleSrc5.Copy()
$ws.Cells.Item(146, 10).PasteSpecial(-4122)
$ws.Cells.Item(146, 11).Value = '5 havo 1'
$styleSrc6.Copy()
$ws.Cells.Item(146, 11).PasteSpecial(-4122)
$ws.Cells.Item(146, 12).Value = 'H5'
$styleSrc6.Copy()
$ws.Cells.Item(146, 12).PasteSpecial(-4122)
$styleSrc5.Copy()
$ws.Cells.Item(146, 14).PasteSpecial(-4122)
$ws.Rows.Item(146).RowHeight = 51

# Row 147
$ws.Cells.Item(147, 2).Value = 18
$ws.Cells.Item(147, 3).Value = 'Bereken'
$ws.Cells.Item(147, 4).Value = 'B2: Algebra'
$styleSrc3.Copy()
$ws.Cells.Item(147, 4).PasteSpecial(-4122)
$ws.Cells.Item(147, 5).Value = 'Substitutie
Grafische rekenmachine
Afronden
Vergelijking'
$styleSrc3.Copy()
$ws.Cells.Item(147, 5).PasteSpecial(-4122)
$ws.Cells.Item(147, 6).Value = 'Twee formules combineren, waarde invullen, vergelijking oplossen (GR), conclusie trekken'
$styleSrc4.Copy()
$ws.Cells.Item(147, 6).PasteSpecial(-4122)
$ws.Cells.Item(147, 7).Value = 'Deel 3 H11'
$ws.Cells.Item(147, 8).Value = '12.3'
$ws.Cells.Item(147, 9).Value = 'Formules en variabelen'
$styleSrc5.Copy()
$ws.Cells.Item(147, 9).PasteSpecial(-4122)
$ws.Cells.Item(147, 10).Value = 'Allerlei formules'
$styleSrc5.Copy()
$ws.Cells.Item(147, 10).PasteSpecial(-4122)
$ws.Cells.Item(147, 11).Value = '5 havo 3'
$styleSrc6.Copy()
$ws.Cells.Item(147, 11).PasteSpecial(-4122)
$ws.Cells.Item(147, 12).Value = 'H2'
$styleSrc6.Copy()
$ws.Cells.Item(147, 12).PasteSpecial(-4122)
$styleSrc5.Copy()
$ws.Cells.Item(147, 14).PasteSpecial(-4122)
$ws.Rows.Item(147).RowHeight = 68

# Row 148
$ws.Cells.Item(148, 2).Value = 19
$ws.Cells.Item(148, 3).Value = 'Toon aan'
$styleSrc3.Copy()
$ws.Cells.Item(148, 3).PasteSpecial(-4122)
$ws.Cells.Item(148, 4).Value = 'C5: Exponentiële verbanden'
$styleSrc3.Copy()
$ws.Cells.Item(148, 4).PasteSpecial(-4122)
$ws.Cells.Item(148, 5).Value = 'Groeifactor
Exponentieel verband'
$styleSrc3.Copy()
$ws.Cells.Item(148, 5).PasteSpecial(-4122)
$ws.Cells.Item(148, 6).Value = 'Exponentieel verband aantonen door groeifactoren te vergelijken'
$styleSrc4.Copy()
$ws.Cells.Item(148, 6).PasteSpecial(-4122)
$ws.Cells.Item(148, 7).Value = 'Deel 3 H9'
$styleSrc3.Copy()
$ws.Cells.Item(148, 7).PasteSpecial(-4122)
$ws.Cells.Item(148, 8).Value = '12.4'
$styleSrc3.Copy()
$ws.Cells.Item(148, 8).PasteSpecial(-4122)
$ws.Cells.Item(148, 9).Value = 'Formules en variabelen'
$styleSrc5.Copy()
$ws.Cells.Item(148, 9).PasteSpecial(-4122)
$ws.Cells.Item(148, 10).Value = 'Lineaire en exponentiele formules'
$styleSrc5.Copy()
$ws.Cells.Item(148, 10).PasteSpecial(-4122)
$ws.Cells.Item(148, 11).Value = '5 havo 1'
$styleSrc6.Copy()
$ws.Cells.Item(148, 11).PasteSpecial(-4122)
$ws.Cells.Item(148, 12).Value = 'H5'
$styleSrc6.Copy()
$ws.Cells.Item(148, 12).PasteSpecial(-4122)
$styleSrc5.Copy()
$ws.Cells.Item(148, 14).PasteSpecial(-4122)
$ws.Rows.Item(148).RowHeight = 48

# Row 149
$ws.Cells.Item(149, 2).Value = 20
$ws.Cells.Item(149, 3).Value = 'Bereken'
$ws.Cells.Item(149, 4).Value = 'C4: Lineaire verbanden'
$styleSrc3.Copy()
$ws.Cells.Item(149, 4).PasteSpecial(-4122)
$ws.Cells.Item(149, 5).Value = 'Lineair verband
Lineair inter-/extrapoleren'
$styleSrc3.Copy()
$ws.Cells.Item(149, 5).PasteSpecial(-4122)
$ws.Cells.Item(149, 6).Value = 'Stapgrootte bepalen, extrapoleren'
$styleSrc4.Copy()
$ws.Cells.Item(149, 6).PasteSpecial(-4122)
$ws.Cells.Item(149, 7).Value = 'Deel 2 H5'
$ws.Cells.Item(149, 8).Value = '12.2'
$ws.Cells.Item(149, 9).Value = 'Lineaire verbanden'
$styleSrc5.Copy()
$ws.Cells.Item(149, 9).PasteSpecial(-4122)
$ws.Cells.Item(149, 10).Value = 'Lineaire en exponentiele formules'
$styleSrc5.Copy()
$ws.Cells.Item(149, 10).PasteSpecial(-4122)
$ws.Cells.Item(149, 11).Value = '5 havo 1'
$styleSrc6.Copy()
$ws.Cells.Item(149, 11).PasteSpecial(-4122)
$ws.Cells.Item(149, 12).Value = 'H1'
$styleSrc6.Copy()
$ws.Cells.Item(149, 12).PasteSpecial(-4122)
$styleSrc5.Copy()
$ws.Cells.Item(149, 14).PasteSpecial(-4122)
$ws.Rows.Item(149).RowHeight = 34

# Row 150
$ws.Cells.Item(150, 2).Value = 21
$ws.Cells.Item(150, 3).Value = 'Bereken'
$styleSrc3.Copy()
$ws.Cells.Item(150, 3).PasteSpecial(-4122)
$ws.Cells.Item(150, 4).Value = 'C4: Lineaire verbanden'
$styleSrc3.Copy()
$ws.Cells.Item(150, 4).PasteSpecial(-4122)
$ws.Cells.Item(150, 5).Value = 'Lineaire formule opstellen
Lineair verband
Richtingscoëfficiënt'
$styleSrc3.Copy()
$ws.Cells.Item(150, 5).PasteSpecial(-4122)
$ws.Cells.Item(150, 6).Value = 'Parameters in lineaire formule bepalen'
$styleSrc4.Copy()
$ws.Cells.Item(150, 6).PasteSpecial(-4122)
$ws.Cells.Item(150, 7).Value = 'Deel 2 H5'
$styleSrc3.Copy()
$ws.Cells.Item(150, 7).PasteSpecial(-4122)
$ws.Cells.Item(150, 8).Value = '12.2'
$styleSrc3.Copy()
$ws.Cells.Item(150, 8).PasteSpecial(-4122)
$ws.Cells.Item(150, 9).Value = 'Lineaire verbanden'
$styleSrc5.Copy()
$ws.Cells.Item(150, 9).PasteSpecial(-4122)
$ws.Cells.Item(150, 10).Value = 'Lineaire en exponentiele formules'
$styleSrc5.Copy()
$ws.Cells.Item(150, 10).PasteSpecial(-4122)
$ws.Cells.Item(150, 11).Value = '5 havo 1'
$styleSrc6.Copy()
$ws.Cells.Item(150, 11).PasteSpecial(-4122)
$ws.Cells.Item(150, 12).Value = 'H4'
$styleSrc6.Copy()
$ws.Cells.Item(150, 12).PasteSpecial(-4122)
$styleSrc5.Copy()
$ws.Cells.Item(150, 14).PasteSpecial(-4122)
$ws.Rows.Item(150).RowHeight = 51

# Row 151
$ws.Cells.Item(151, 1).Value = 'Bioscoopbezoek'
$ws.Cells.Item(151, 2).Value = 22
$ws.Cells.Item(151, 3).Value = 'Onderzoek'
$ws.Cells.Item(151, 4).Value = 'C1: Tabellen
C2: Grafieken'
$styleSrc3.Copy()
$ws.Cells.Item(151, 4).PasteSpecial(-4122)
$ws.Cells.Item(151, 5).Value = 'Gemiddelde
Meer dan twee variabelen
Tabel
Staafdiagram/histogram
Statistische visualisatie
Lijndiagram/frequentiepolygoon
Redeneren'
$styleSrc3.Copy()
$ws.Cells.Item(151, 5).PasteSpecial(-4122)
$ws.Cells.Item(151, 6).Value = 'Juiste gegeven uit tabel halen, Gegevens verwerken uit tekst, lijndiagram en staafdiagram, rekenen met gemiddelde'
$styleSrc4.Copy()
$ws.Cells.Item(151, 6).PasteSpecial(-4122)
$ws.Cells.Item(151, 7).Value = 'Deel 1 H3'
$ws.Cells.Item(151, 8).Value = '12.1'
$ws.Cells.Item(151, 9).Value = 'Tabellen en grafieken'
$styleSrc5.Copy()
$ws.Cells.Item(151, 9).PasteSpecial(-4122)
$ws.Cells.Item(151, 10).Value = 'Conclusies uit data'
$styleSrc5.Copy()
$ws.Cells.Item(151, 10).PasteSpecial(-4122)
$ws.Cells.Item(151, 11).Value = '5 havo 5'
$styleSrc7.Copy()
$ws.Cells.Item(151, 11).PasteSpecial(-4122)
$ws.Cells.Item(151, 12).Value = 'H2'
$styleSrc7.Copy()
$ws.Cells.Item(151, 12).PasteSpecial(-4122)
$styleSrc5.Copy()
$ws.Cells.Item(151, 14).PasteSpecial(-4122)
$ws.Rows.Item(151).RowHeight = 136

$ws.Application.CutCopyMode = $false

# Refresh the view to match the post-edit selection/scroll position
$ws.Application.Goto($ws.Range("A128"), $true)
$ws.Range("N130:N151").Select()
